{"js": "// \"Version 2.\" -> \"Version 1.\" while preserving the existing run/proofErr/\n// bookmark structure as closely as possible (matches the target OOXML diff):\n//   - merge the \"Versi\" + \"on\" runs into a single \"Version\" run\n//   - change the \" 2\" run's text to \" 1.\"\n//   - remove the trailing \".\" run entirely\n\nconst body = context.document.body;\n\n// 1) Merge \"Versi\" + \"on\" (the word spanned by the spellStart/spellEnd\n//    proofErr) into a single \"Version\" run. Replacing the combined range in\n//    one insertText call forces the two backing runs to merge into one.\nconst versiResults = body.search(\"Versi\", { matchCase: true });\nversiResults.load(\"items\");\nconst onResults = body.search(\"on\", { matchCase: true });\nonResults.load(\"items\");\nawait context.sync();\n\nconst versiRange = versiResults.items[0];\nconst onRange = onResults.items[0];\nconst versionRange = versiRange.expandTo(onRange);\nversionRange.insertText(\"Version\", Word.InsertLocation.replace);\nawait context.sync();\n\n// 2) Change \" 2\" to \" 1.\" in place (this run stays right before the bookmark).\nconst spaceTwoResults = body.search(\" 2\", { matchCase: true });\nspaceTwoResults.load(\"items\");\nawait context.sync();\nspaceTwoResults.items[0].insertText(\" 1.\", Word.InsertLocation.replace);\nawait context.sync();\n\n// 3) Delete the now-redundant trailing \".\" run. Anchor on the _GoBack\n//    bookmark (which sits right before it) and the paragraph end (which\n//    sits right after it, before the paragraph mark) rather than searching\n//    for \".\" again, so this does not depend on the text changes above.\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nconst bookmarkRange = body.getBookmarkRangeOrNullObject(\"_GoBack\");\nawait context.sync();\n\nconst paragraph = paragraphs.items[0];\nconst paragraphEnd = paragraph.getRange(\"End\");\nconst trailingDot = bookmarkRange.expandTo(paragraphEnd);\ntrailingDot.delete();\nawait context.sync();\n", "ps1": "# \"Version 2.\" -> \"Version 1.\" while preserving the existing run/proofErr/\n# bookmark structure as closely as possible (matches the target OOXML diff):\n#   - merge the \"Versi\" + \"on\" runs into a single \"Version\" run\n#   - change the \" 2\" run's text to \" 1.\"\n#   - remove the trailing \".\" run entirely\n\n$d = $word.ActiveDocument\n\n# 1) Locate \"Versi\" then the \"on\" that immediately follows it (the word is\n#    spanned by a spellStart/spellEnd proofErr pair and currently split\n#    across two runs). Union the two finds into one range and run a\n#    Find/Replace over it: even though the replacement text is identical to\n#    what's already there, routing the edit through Find/Replace forces the\n#    two backing runs to normalize into a single \"Version\" run.\n$versiRange = $d.Content\n$versiFind = $versiRange.Find\n$versiFind.Text = \"Versi\"\n$versiFind.Execute() | Out-Null\n\n$onRange = $d.Range($versiRange.End, $d.Content.End)\n$onFind = $onRange.Find\n$onFind.Text = \"on\"\n$onFind.Execute() | Out-Null\n\n$wordRange = $d.Range($versiRange.Start, $onRange.End)\n$mergeFind = $wordRange.Find\n$mergeFind.Text = $wordRange.Text\n$mergeFind.Replacement.Text = $wordRange.Text\n$mergeFind.Execute([ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, 2)\n\n# 2) Change \" 2\" to \" 1.\" in place (this run stays right before the bookmark).\n$spaceTwoRange = $d.Content\n$spaceTwoFind = $spaceTwoRange.Find\n$spaceTwoFind.Text = \" 2\"\n$spaceTwoFind.Execute() | Out-Null\n$spaceTwoRange.Text = \" 1.\"\n\n# 3) Delete the now-redundant trailing \".\" run. Anchor on the _GoBack\n#    bookmark (which sits right before it) and the paragraph end (which sits\n#    right after it, before the paragraph mark) so this is robust to the\n#    text shifts caused by steps 1-2.\n$bookmark = $d.Bookmarks(\"_GoBack\")\n$paragraph = $d.Paragraphs.First\n$trailingDot = $d.Range($bookmark.End, $paragraph.Range.End - 1)\n$trailingDot.Text = \"\"\n"}
